$d = $word.ActiveDocument

# Locate the "PS: ..." paragraph that the new paragraph must follow.
$findRange = $d.Content
$found = $findRange.Find.Execute("PS: A quantidade de produtos vendidos foi de 10",
                                  $true, $false, $false, $false, $false,
                                  $true, 1, $false, "", 0)

# Move to the end of that paragraph (after the found text, still inside
# the same paragraph, before its paragraph mark).
$endRange = $d.Range($findRange.End, $findRange.End)

# Insert a new paragraph mark right after it.
$endRange.InsertParagraphAfter()

# The newly created paragraph is now the last one in the document; give it
# the "MacroText" style and the requested text.
$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "testando estilos aleatorios"
$newPara.Style = "MacroText"
